$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 15-47: columns D, J, K, L, M, O, P changed per diff ---
# Row 15
$ws.Cells.Item(15, 4).Value = 44188
$ws.Cells.Item(15, 10).Value = 100
$ws.Cells.Item(15, 11).Value = 42000
$ws.Cells.Item(15, 12).Value = 44000
$ws.Cells.Item(15, 13).Value = 43000
$ws.Cells.Item(15, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(15, 16).Value = 1720

# Row 16
$ws.Cells.Item(16, 4).Value = 44244
$ws.Cells.Item(16, 10).Value = 100
$ws.Cells.Item(16, 11).Value = 25000
$ws.Cells.Item(16, 12).Value = 26000
$ws.Cells.Item(16, 13).Value = 25500
$ws.Cells.Item(16, 15).Value = 'Región del Maule'
$ws.Cells.Item(16, 16).Value = 1020

# Row 17
$ws.Cells.Item(17, 4).Value = 44230
$ws.Cells.Item(17, 10).Value = 100
$ws.Cells.Item(17, 11).Value = 35000
$ws.Cells.Item(17, 12).Value = 36000
$ws.Cells.Item(17, 13).Value = 35500
$ws.Cells.Item(17, 15).Value = 'Región del Maule'
$ws.Cells.Item(17, 16).Value = 1420

# Row 18
$ws.Cells.Item(18, 4).Value = 44265
$ws.Cells.Item(18, 10).Value = 100
$ws.Cells.Item(18, 11).Value = 22000
$ws.Cells.Item(18, 12).Value = 24000
$ws.Cells.Item(18, 13).Value = 23000
$ws.Cells.Item(18, 15).Value = 'Región del Maule'
$ws.Cells.Item(18, 16).Value = 920

# Row 19
$ws.Cells.Item(19, 4).Value = 44272
$ws.Cells.Item(19, 10).Value = 100
$ws.Cells.Item(19, 11).Value = 22000
$ws.Cells.Item(19, 12).Value = 24000
$ws.Cells.Item(19, 13).Value = 23000
$ws.Cells.Item(19, 15).Value = 'Región del Maule'
$ws.Cells.Item(19, 16).Value = 920

# Row 20
$ws.Cells.Item(20, 4).Value = 44574
$ws.Cells.Item(20, 10).Value = 100
$ws.Cells.Item(20, 11).Value = 30000
$ws.Cells.Item(20, 12).Value = 32000
$ws.Cells.Item(20, 13).Value = 31000
$ws.Cells.Item(20, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(20, 16).Value = 1240

# Row 21
$ws.Cells.Item(21, 4).Value = 44236
$ws.Cells.Item(21, 10).Value = 100
$ws.Cells.Item(21, 11).Value = 25000
$ws.Cells.Item(21, 12).Value = 26000
$ws.Cells.Item(21, 13).Value = 25500
$ws.Cells.Item(21, 15).Value = 'Región del Maule'
$ws.Cells.Item(21, 16).Value = 1020

# Row 22
$ws.Cells.Item(22, 4).Value = 44587
$ws.Cells.Item(22, 10).Value = 220
$ws.Cells.Item(22, 11).Value = 23000
$ws.Cells.Item(22, 12).Value = 24000
$ws.Cells.Item(22, 13).Value = 23545
$ws.Cells.Item(22, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(22, 16).Value = 942

# Row 23
$ws.Cells.Item(23, 4).Value = 44313
$ws.Cells.Item(23, 10).Value = 100
$ws.Cells.Item(23, 11).Value = 30000
$ws.Cells.Item(23, 12).Value = 32000
$ws.Cells.Item(23, 13).Value = 31000
$ws.Cells.Item(23, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(23, 16).Value = 1240

# Row 24
$ws.Cells.Item(24, 4).Value = 44203
$ws.Cells.Item(24, 10).Value = 100
$ws.Cells.Item(24, 11).Value = 25000
$ws.Cells.Item(24, 12).Value = 26000
$ws.Cells.Item(24, 13).Value = 25500
$ws.Cells.Item(24, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(24, 16).Value = 1020

# Row 25
$ws.Cells.Item(25, 4).Value = 44558
$ws.Cells.Item(25, 10).Value = 250
$ws.Cells.Item(25, 11).Value = 15000
$ws.Cells.Item(25, 12).Value = 16000
$ws.Cells.Item(25, 13).Value = 15400
$ws.Cells.Item(25, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(25, 16).Value = 616

# Row 26
$ws.Cells.Item(26, 4).Value = 44210
$ws.Cells.Item(26, 10).Value = 100
$ws.Cells.Item(26, 11).Value = 32000
$ws.Cells.Item(26, 12).Value = 34000
$ws.Cells.Item(26, 13).Value = 33000
$ws.Cells.Item(26, 15).Value = 'Región del Maule'
$ws.Cells.Item(26, 16).Value = 1320

# Row 27
$ws.Cells.Item(27, 4).Value = 44316
$ws.Cells.Item(27, 10).Value = 100
$ws.Cells.Item(27, 11).Value = 26000
$ws.Cells.Item(27, 12).Value = 27000
$ws.Cells.Item(27, 13).Value = 26500
$ws.Cells.Item(27, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(27, 16).Value = 1060

# Row 28
$ws.Cells.Item(28, 4).Value = 44692
$ws.Cells.Item(28, 10).Value = 100
$ws.Cells.Item(28, 11).Value = 25000
$ws.Cells.Item(28, 12).Value = 26000
$ws.Cells.Item(28, 13).Value = 25500
$ws.Cells.Item(28, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(28, 16).Value = 1020

# Row 29
$ws.Cells.Item(29, 4).Value = 44923
$ws.Cells.Item(29, 10).Value = 100
$ws.Cells.Item(29, 11).Value = 35000
$ws.Cells.Item(29, 12).Value = 36000
$ws.Cells.Item(29, 13).Value = 35500
$ws.Cells.Item(29, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(29, 16).Value = 1420

# Row 30
$ws.Cells.Item(30, 4).Value = 44923
$ws.Cells.Item(30, 10).Value = 100
$ws.Cells.Item(30, 11).Value = 35000
$ws.Cells.Item(30, 12).Value = 36000
$ws.Cells.Item(30, 13).Value = 35500
$ws.Cells.Item(30, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(30, 16).Value = 1420

# Row 31
$ws.Cells.Item(31, 4).Value = 44568
$ws.Cells.Item(31, 10).Value = 200
$ws.Cells.Item(31, 11).Value = 25000
$ws.Cells.Item(31, 12).Value = 26000
$ws.Cells.Item(31, 13).Value = 25500
$ws.Cells.Item(31, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(31, 16).Value = 1020

# Row 32
$ws.Cells.Item(32, 4).Value = 44308
$ws.Cells.Item(32, 10).Value = 100
$ws.Cells.Item(32, 11).Value = 28000
$ws.Cells.Item(32, 12).Value = 30000
$ws.Cells.Item(32, 13).Value = 29000
$ws.Cells.Item(32, 15).Value = 'Región del Maule'
$ws.Cells.Item(32, 16).Value = 1160

# Row 33
$ws.Cells.Item(33, 4).Value = 44624
$ws.Cells.Item(33, 10).Value = 150
$ws.Cells.Item(33, 11).Value = 25000
$ws.Cells.Item(33, 12).Value = 26000
$ws.Cells.Item(33, 13).Value = 25467
$ws.Cells.Item(33, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(33, 16).Value = 1019

# Row 34
$ws.Cells.Item(34, 4).Value = 44897
$ws.Cells.Item(34, 10).Value = 100
$ws.Cells.Item(34, 11).Value = 38000
$ws.Cells.Item(34, 12).Value = 40000
$ws.Cells.Item(34, 13).Value = 39000
$ws.Cells.Item(34, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(34, 16).Value = 1560

# Row 35
$ws.Cells.Item(35, 4).Value = 44603
$ws.Cells.Item(35, 10).Value = 130
$ws.Cells.Item(35, 11).Value = 25000
$ws.Cells.Item(35, 12).Value = 26000
$ws.Cells.Item(35, 13).Value = 25385
$ws.Cells.Item(35, 15).Value = 'Región del Maule'
$ws.Cells.Item(35, 16).Value = 1015

# Row 36
$ws.Cells.Item(36, 4).Value = 44602
$ws.Cells.Item(36, 10).Value = 130
$ws.Cells.Item(36, 11).Value = 20000
$ws.Cells.Item(36, 12).Value = 21000
$ws.Cells.Item(36, 13).Value = 20385
$ws.Cells.Item(36, 15).Value = 'Región del Maule'
$ws.Cells.Item(36, 16).Value = 815

# Row 37
$ws.Cells.Item(37, 4).Value = 44560
$ws.Cells.Item(37, 10).Value = 100
$ws.Cells.Item(37, 11).Value = 25000
$ws.Cells.Item(37, 12).Value = 26000
$ws.Cells.Item(37, 13).Value = 25500
$ws.Cells.Item(37, 15).Value = 'Región del Maule'
$ws.Cells.Item(37, 16).Value = 1020

# Row 38
$ws.Cells.Item(38, 4).Value = 44636
$ws.Cells.Item(38, 10).Value = 80
$ws.Cells.Item(38, 11).Value = 22000
$ws.Cells.Item(38, 12).Value = 23000
$ws.Cells.Item(38, 13).Value = 22375
$ws.Cells.Item(38, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(38, 16).Value = 895

# Row 39
$ws.Cells.Item(39, 4).Value = 44609
$ws.Cells.Item(39, 10).Value = 200
$ws.Cells.Item(39, 11).Value = 26000
$ws.Cells.Item(39, 12).Value = 28000
$ws.Cells.Item(39, 13).Value = 27000
$ws.Cells.Item(39, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(39, 16).Value = 1080

# Row 40
$ws.Cells.Item(40, 4).Value = 44328
$ws.Cells.Item(40, 10).Value = 100
$ws.Cells.Item(40, 11).Value = 32000
$ws.Cells.Item(40, 12).Value = 34000
$ws.Cells.Item(40, 13).Value = 33000
$ws.Cells.Item(40, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(40, 16).Value = 1320

# Row 41
$ws.Cells.Item(41, 4).Value = 44657
$ws.Cells.Item(41, 10).Value = 140
$ws.Cells.Item(41, 11).Value = 18000
$ws.Cells.Item(41, 12).Value = 20000
$ws.Cells.Item(41, 13).Value = 19143
$ws.Cells.Item(41, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(41, 16).Value = 766

# Row 42
$ws.Cells.Item(42, 4).Value = 44685
$ws.Cells.Item(42, 10).Value = 150
$ws.Cells.Item(42, 11).Value = 25000
$ws.Cells.Item(42, 12).Value = 26000
$ws.Cells.Item(42, 13).Value = 25467
$ws.Cells.Item(42, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(42, 16).Value = 1019

# Row 43
$ws.Cells.Item(43, 4).Value = 44342
$ws.Cells.Item(43, 10).Value = 100
$ws.Cells.Item(43, 11).Value = 28000
$ws.Cells.Item(43, 12).Value = 30000
$ws.Cells.Item(43, 13).Value = 29000
$ws.Cells.Item(43, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(43, 16).Value = 1160

# Row 44
$ws.Cells.Item(44, 4).Value = 44194
$ws.Cells.Item(44, 10).Value = 100
$ws.Cells.Item(44, 11).Value = 30000
$ws.Cells.Item(44, 12).Value = 32000
$ws.Cells.Item(44, 13).Value = 31000
$ws.Cells.Item(44, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(44, 16).Value = 1240

# Row 45
$ws.Cells.Item(45, 4).Value = 44638
$ws.Cells.Item(45, 10).Value = 220
$ws.Cells.Item(45, 11).Value = 20000
$ws.Cells.Item(45, 12).Value = 22000
$ws.Cells.Item(45, 13).Value = 21091
$ws.Cells.Item(45, 15).Value = 'Región del Maule'
$ws.Cells.Item(45, 16).Value = 844

# Row 46
$ws.Cells.Item(46, 4).Value = 44651
$ws.Cells.Item(46, 10).Value = 250
$ws.Cells.Item(46, 11).Value = 28000
$ws.Cells.Item(46, 12).Value = 30000
$ws.Cells.Item(46, 13).Value = 28960
$ws.Cells.Item(46, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(46, 16).Value = 1158

# Row 47
$ws.Cells.Item(47, 4).Value = 44216
$ws.Cells.Item(47, 10).Value = 100
$ws.Cells.Item(47, 11).Value = 26000
$ws.Cells.Item(47, 12).Value = 28000
$ws.Cells.Item(47, 13).Value = 27000
$ws.Cells.Item(47, 15).Value = 'Región del Maule'
$ws.Cells.Item(47, 16).Value = 1080

# --- Append new rows 48-50 (full rows), matching the layout/format of prior rows ---
# Row 48
$ws.Cells.Item(48, 1).Value = 11
$ws.Cells.Item(48, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(48, 3).Value = 'Bíobío'
$ws.Cells.Item(48, 4).Value = 44629
$ws.Cells.Item(48, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(48, 5).Value = 8
$ws.Cells.Item(48, 6).Value = 100112030
$ws.Cells.Item(48, 7).Value = 'Poroto granado'
$ws.Cells.Item(48, 8).Value = 'Sin especificar'
$ws.Cells.Item(48, 9).Value = 'Primera'
$ws.Cells.Item(48, 10).Value = 110
$ws.Cells.Item(48, 11).Value = 25000
$ws.Cells.Item(48, 12).Value = 26000
$ws.Cells.Item(48, 13).Value = 25455
$ws.Cells.Item(48, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(48, 15).Value = 'Región del Maule'
$ws.Cells.Item(48, 16).Value = 1018
$ws.Cells.Item(48, 17).Value = 25
$ws.Cells.Item(48, 18).Value = 'Hortaliza'

# Row 49
$ws.Cells.Item(49, 1).Value = 11
$ws.Cells.Item(49, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(49, 3).Value = 'Bíobío'
$ws.Cells.Item(49, 4).Value = 44706
$ws.Cells.Item(49, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(49, 5).Value = 8
$ws.Cells.Item(49, 6).Value = 100112030
$ws.Cells.Item(49, 7).Value = 'Poroto granado'
$ws.Cells.Item(49, 8).Value = 'Sin especificar'
$ws.Cells.Item(49, 9).Value = 'Primera'
$ws.Cells.Item(49, 10).Value = 100
$ws.Cells.Item(49, 11).Value = 23000
$ws.Cells.Item(49, 12).Value = 25000
$ws.Cells.Item(49, 13).Value = 24000
$ws.Cells.Item(49, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(49, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(49, 16).Value = 960
$ws.Cells.Item(49, 17).Value = 25
$ws.Cells.Item(49, 18).Value = 'Hortaliza'

# Row 50
$ws.Cells.Item(50, 1).Value = 11
$ws.Cells.Item(50, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(50, 3).Value = 'Bíobío'
$ws.Cells.Item(50, 4).Value = 44580
$ws.Cells.Item(50, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(50, 5).Value = 8
$ws.Cells.Item(50, 6).Value = 100112030
$ws.Cells.Item(50, 7).Value = 'Poroto granado'
$ws.Cells.Item(50, 8).Value = 'Sin especificar'
$ws.Cells.Item(50, 9).Value = 'Primera'
$ws.Cells.Item(50, 10).Value = 100
$ws.Cells.Item(50, 11).Value = 28000
$ws.Cells.Item(50, 12).Value = 30000
$ws.Cells.Item(50, 13).Value = 29000
$ws.Cells.Item(50, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(50, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(50, 16).Value = 1160
$ws.Cells.Item(50, 17).Value = 25
$ws.Cells.Item(50, 18).Value = 'Hortaliza'

